# Auto-generated Excel COM-interop script
# Applies scheduled market-price refresh updates to the Leve profit tables
# across all job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
#
# For each affected row, columns H-N (currentAveragePrice, currentAveragePriceNQ,
# currentAveragePriceHQ, LevePriceNQ, LevePriceHQ, LeveProfitNQ, LeveProfitHQ) are
# refreshed with newly retrieved market data. Some cells go from a N/A (blank)
# state to a computed value, and vice versa, depending on market data availability;
# those are handled with ClearContents() / Value assignment as appropriate.

$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")

$updates = @{
    "H17" = 877.6923
    "J17" = 877.6923
    "L17" = 2633.0769
    "N17" = -2969.0769
    "H18" = 725.25
    "I18" = 725.25
    "K18" = 725.25
    "M18" = -441.25
    "H40" = 2421.1
    "I40" = 988.4
    "J40" = 2898.6667
    "K40" = 988.4
    "L40" = 2898.6667
    "M40" = -813.4
    "N40" = -3248.6667
    "H62" = 9644.909
    "J62" = 12071
    "L62" = 12071
    "N62" = -13319
    "H65" = 9644.909
    "J65" = 12071
    "L65" = 60355
    "N65" = -66595
    "H88" = 11090.272
    "J88" = 11199.4
    "L88" = 11199.4
    "N88" = -12011.4
    "H91" = 11090.272
    "J91" = 11199.4
    "L91" = 11199.4
    "N91" = -14007.4
    "H97" = 3869.375
    "J97" = 3869.375
    "L97" = 11608.125
    "N97" = -12600.125
    "H112" = 1626.05
    "I112" = 1494
    "J112" = 1682.6428
    "K112" = 4482
    "L112" = 5047.928400000001
    "M112" = -3374
    "N112" = -7263.928400000001
    "H126" = 0
    "J126" = 0
    "L126" = 0
    "H129" = 1155.5
    "J129" = 1712.5714
    "L129" = 5137.7142
    "N129" = -15137.7142
    "H138" = 3647.4143
    "J138" = 3718.356
    "L138" = 11155.068
    "N138" = -21435.068
    "H141" = 3100.5833
    "I141" = 2800.889
    "K141" = 8402.667000000001
    "M141" = -3222.667000000001
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

foreach ($ref in @("N126")) {
    $ws.Range($ref).ClearContents()
}

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")

$updates = @{
    "H44" = 62497.5
    "J44" = 62497.5
    "L44" = 62497.5
    "N44" = -63473.5
    "H55" = 35000
    "I55" = 20000
    "J55" = 50000
    "K55" = 20000
    "L55" = 50000
    "M55" = -19685
    "N55" = -50630
    "H63" = 11287.833
    "I63" = 13489.8
    "J63" = 9715
    "K63" = 13489.8
    "L63" = 9715
    "M63" = -12803.8
    "N63" = -11087
    "H66" = 11287.833
    "I66" = 13489.8
    "J66" = 9715
    "K66" = 67449
    "L66" = 48575
    "M66" = -64017
    "N66" = -55439
    "H74" = 4945.826
    "I74" = 4309.294
    "K74" = 4309.294
    "M74" = -3435.294
    "H77" = 4945.826
    "I77" = 4309.294
    "K77" = 21546.47
    "M77" = -17178.47
    "H80" = 0
    "I80" = 0
    "K80" = 0
    "H83" = 0
    "I83" = 0
    "K83" = 0
    "H92" = 49999
    "J92" = 49999
    "L92" = 49999
    "N92" = -54991
    "H97" = 818.1
    "I97" = 135.125
    "K97" = 135.125
    "M97" = 360.875
    "H102" = 4387.067
    "I102" = 2983.8333
    "K102" = 2983.8333
    "M102" = -1361.8333
    "H122" = 3833050.8
    "I122" = 4446075
    "K122" = 13338225
    "M122" = -13335775
    "H132" = 2735.5
    "I132" = 2645.6155
    "J132" = 3125
    "K132" = 7936.8465
    "L132" = 9375
    "M132" = -5406.8465
    "N132" = -14435
    "H138" = 119998.5
    "J138" = 119998.5
    "L138" = 119998.5
    "N138" = -130278.5
    "H139" = 62905
    "J139" = 62905
    "L139" = 62905
    "N139" = -73185
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

foreach ($ref in @("M80","M83")) {
    $ws.Range($ref).ClearContents()
}

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")

$updates = @{
    "H82" = 34513.5
    "J82" = 58999
    "L82" = 58999
    "N82" = -59765
    "H85" = 34513.5
    "J85" = 58999
    "L85" = 58999
    "N85" = -61651
    "H94" = 151736.52
    "I94" = 217830.94
    "J94" = 663.5714
    "K94" = 217830.94
    "L94" = 663.5714
    "M94" = -217379.94
    "N94" = -1565.5714
    "H99" = 5064.7144
    "I99" = 4900.7
    "J99" = 5474.75
    "K99" = 4900.7
    "L99" = 5474.75
    "M99" = -3402.7
    "N99" = -8470.75
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")

$updates = @{
    "H31" = 3771.1765
    "I31" = 1565.1428
    "K31" = 1565.1428
    "M31" = -1270.1428
    "H34" = 3771.1765
    "I34" = 1565.1428
    "K34" = 1565.1428
    "M34" = -1363.1428
    "H58" = 2543.5
    "I58" = 3087.5
    "J58" = 1999.5
    "K58" = 3087.5
    "L58" = 1999.5
    "M58" = -2884.5
    "N58" = -2405.5
    "H136" = 2543.5
    "I136" = 3087.5
    "J136" = 1999.5
    "K136" = 9262.5
    "L136" = 5998.5
    "M136" = -6712.5
    "N136" = -11098.5
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")

$updates = @{
    "H5" = 1541.826
    "I5" = 132.2
    "J5" = 2626.1538
    "K5" = 396.6
    "L5" = 7878.4614
    "M5" = -284.6
    "N5" = -8102.4614
    "H107" = 2598.2144
    "I107" = 2175.625
    "J107" = 3161.6667
    "K107" = 6526.875
    "L107" = 9485.000100000001
    "M107" = -4606.875
    "N107" = -13325.0001
    "H131" = 3385810.5
    "I131" = 250480
    "J131" = 4169643
    "K131" = 751440
    "L131" = 12508929
    "M131" = -746400
    "N131" = -12519009
    "H132" = 1100
    "I132" = 1100
    "J132" = 0
    "K132" = 9900
    "L132" = 0
    "M132" = -7370
    "H135" = 1541.826
    "I135" = 132.2
    "J135" = 2626.1538
    "K135" = 1189.8
    "L135" = 23635.3842
    "M135" = 1345.2
    "N135" = -28705.3842
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

foreach ($ref in @("N132")) {
    $ws.Range($ref).ClearContents()
}

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")

$updates = @{
    "H97" = 4720.3706
    "I97" = 1189.7273
    "K97" = 1189.7273
    "M97" = -693.7273
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")

$updates = @{
    "H93" = 4468.4375
    "I93" = 2249.7144
    "K93" = 2249.7144
    "M93" = -1001.7144
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")

$updates = @{
    "H96" = 4540.636
    "J96" = 4518.625
    "L96" = 4518.625
    "N96" = -7264.625
    "H126" = 2974.75
    "I126" = 2966.3333
    "K126" = 8898.999899999999
    "M126" = -6428.999899999999
    "H136" = 3558.25
    "I136" = 3355.4443
    "J136" = 4166.6665
    "K136" = 10066.3329
    "L136" = 12499.9995
    "M136" = -7516.332900000001
    "N136" = -17599.9995
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
